# Rename the header row column names so the "old"/"new" comparison-suffix
# naming scheme becomes the explicit format-version naming scheme
# (<suffix>_old -> <suffix>_FV2310, <suffix>_new -> <suffix>_FV2404).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Freeze the header row (split below row 1, keep it visible while scrolling).
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the used range into a native Excel table ("Table1") spanning the full
# data extent, with the renamed headers becoming the table's column headers.
$dataRange = $ws.Range("A1:U64")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
